$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2024-08-21 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-22 Thursday", 2)

# Update each table cell value (positional, row-major, to handle duplicate source text)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "10+55="
$t.Cell(1, 2).Range.Text = "35-22="
$t.Cell(1, 3).Range.Text = "43+15="
$t.Cell(1, 4).Range.Text = "93-4="
$t.Cell(1, 5).Range.Text = "62+1="
$t.Cell(2, 1).Range.Text = "3+17="
$t.Cell(2, 2).Range.Text = "67-10="
$t.Cell(2, 3).Range.Text = "25+35="
$t.Cell(2, 4).Range.Text = "10+1="
$t.Cell(2, 5).Range.Text = "11+86="
$t.Cell(3, 1).Range.Text = "0+34="
$t.Cell(3, 2).Range.Text = "94-67="
$t.Cell(3, 3).Range.Text = "40+47="
$t.Cell(3, 4).Range.Text = "27-20="
$t.Cell(3, 5).Range.Text = "6+36="
$t.Cell(4, 1).Range.Text = "64+7="
$t.Cell(4, 2).Range.Text = "59-51="
$t.Cell(4, 3).Range.Text = "84-68="
$t.Cell(4, 4).Range.Text = "75-27="
$t.Cell(4, 5).Range.Text = "13-12="
$t.Cell(5, 1).Range.Text = "51+13="
$t.Cell(5, 2).Range.Text = "87-66="
$t.Cell(5, 3).Range.Text = "41+47="
$t.Cell(5, 4).Range.Text = "39+47="
$t.Cell(5, 5).Range.Text = "23+7="
$t.Cell(6, 1).Range.Text = "41+44="
$t.Cell(6, 2).Range.Text = "0+27="
$t.Cell(6, 3).Range.Text = "44-32="
$t.Cell(6, 4).Range.Text = "99-43="
$t.Cell(6, 5).Range.Text = "26+60="
$t.Cell(7, 1).Range.Text = "92-44="
$t.Cell(7, 2).Range.Text = "41+58="
$t.Cell(7, 3).Range.Text = "96-70="
$t.Cell(7, 4).Range.Text = "11+74="
$t.Cell(7, 5).Range.Text = "54-25="
$t.Cell(8, 1).Range.Text = "33+63="
$t.Cell(8, 2).Range.Text = "13+36="
$t.Cell(8, 3).Range.Text = "70+15="
$t.Cell(8, 4).Range.Text = "80-27="
$t.Cell(8, 5).Range.Text = "48-5="
$t.Cell(9, 1).Range.Text = "98-16="
$t.Cell(9, 2).Range.Text = "98-31="
$t.Cell(9, 3).Range.Text = "50-41="
$t.Cell(9, 4).Range.Text = "76-67="
$t.Cell(9, 5).Range.Text = "63+12="
$t.Cell(10, 1).Range.Text = "19+28="
$t.Cell(10, 2).Range.Text = "1+22="
$t.Cell(10, 3).Range.Text = "88-29="
$t.Cell(10, 4).Range.Text = "88-67="
$t.Cell(10, 5).Range.Text = "25-17="
$t.Cell(11, 1).Range.Text = "40+54="
$t.Cell(11, 2).Range.Text = "16+29="
$t.Cell(11, 3).Range.Text = "39+20="
$t.Cell(11, 4).Range.Text = "2+16="
$t.Cell(11, 5).Range.Text = "34+59="
$t.Cell(12, 1).Range.Text = "31+53="
$t.Cell(12, 2).Range.Text = "5+83="
$t.Cell(12, 3).Range.Text = "53+5="
$t.Cell(12, 4).Range.Text = "93+2="
$t.Cell(12, 5).Range.Text = "70+0="
$t.Cell(13, 1).Range.Text = "66-57="
$t.Cell(13, 2).Range.Text = "57+16="
$t.Cell(13, 3).Range.Text = "36+29="
$t.Cell(13, 4).Range.Text = "73+7="
$t.Cell(13, 5).Range.Text = "28+63="
$t.Cell(14, 1).Range.Text = "11+53="
$t.Cell(14, 2).Range.Text = "80-41="
$t.Cell(14, 3).Range.Text = "18-11="
$t.Cell(14, 4).Range.Text = "52-17="
$t.Cell(14, 5).Range.Text = "99-7="
$t.Cell(15, 1).Range.Text = "41+8="
$t.Cell(15, 2).Range.Text = "63-2="
$t.Cell(15, 3).Range.Text = "31+59="
$t.Cell(15, 4).Range.Text = "98-9="
$t.Cell(15, 5).Range.Text = "43+1="
$t.Cell(16, 1).Range.Text = "86-52="
$t.Cell(16, 2).Range.Text = "33-17="
$t.Cell(16, 3).Range.Text = "86-0="
$t.Cell(16, 4).Range.Text = "18+80="
$t.Cell(16, 5).Range.Text = "44+0="
$t.Cell(17, 1).Range.Text = "81-55="
$t.Cell(17, 2).Range.Text = "17+1="
$t.Cell(17, 3).Range.Text = "62-44="
$t.Cell(17, 4).Range.Text = "36+3="
$t.Cell(17, 5).Range.Text = "25+53="
$t.Cell(18, 1).Range.Text = "65+23="
$t.Cell(18, 2).Range.Text = "65+26="
$t.Cell(18, 3).Range.Text = "68-6="
$t.Cell(18, 4).Range.Text = "80-68="
$t.Cell(18, 5).Range.Text = "56+36="
$t.Cell(19, 1).Range.Text = "5-2="
$t.Cell(19, 2).Range.Text = "87-21="
$t.Cell(19, 3).Range.Text = "21+28="
$t.Cell(19, 4).Range.Text = "51-3="
$t.Cell(19, 5).Range.Text = "54-29="
$t.Cell(20, 1).Range.Text = "61+29="
$t.Cell(20, 2).Range.Text = "54+4="
$t.Cell(20, 3).Range.Text = "81+17="
$t.Cell(20, 4).Range.Text = "91-91="
$t.Cell(20, 5).Range.Text = "14+78="
